$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before the old "delmar" column (AE) for the new
# "coronado" city column. This shifts every column from AE onward one to
# the right (delmar -> AF, elcajon -> AG, ... springvalley -> BA). ---
$ws.Range("AE1").EntireColumn.Insert()
$ws.Range("AE1").Value = "coronado"

# --- Extend the (hidden) AutoFilter database defined name by one column
# to account for the newly inserted column. ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$AY`$17"
    }
}

# --- Append the new day's row (27 March 2020 = serial 43917) of data.
# Copy the prior row's formatting first so the date cell keeps the same
# date number-format style as the rest of column A. ---
$ws.Range("A22").Copy($ws.Range("A23")) | Out-Null

$ws.Range("A23").Value = 43917
$ws.Range("B23").Value = 417
$ws.Range("C23").Value = 3
$ws.Range("E23").Value = 4
$ws.Range("G23").Value = 84
$ws.Range("I23").Value = 106
$ws.Range("K23").Value = 73
$ws.Range("M23").Value = 65
$ws.Range("O23").Value = 31
$ws.Range("Q23").Value = 29
$ws.Range("S23").Value = 20
$ws.Range("U23").Value = 2
$ws.Range("W23").Value = 166
$ws.Range("X23").Value = 249
$ws.Range("Y23").Value = 2
$ws.Range("Z23").Value = 85
$ws.Range("AA23").Value = 38
$ws.Range("AB23").Value = 5
$ws.Range("AC23").Value = 18
$ws.Range("AD23").Value = 21
$ws.Range("AE23").Value = 1
$ws.Range("AF23").Value = 6
$ws.Range("AG23").Value = 20
$ws.Range("AH23").Value = 14
$ws.Range("AI23").Value = 7
$ws.Range("AJ23").Value = 3
$ws.Range("AK23").Value = 1
$ws.Range("AL23").Value = 6
$ws.Range("AM23").Value = 8
$ws.Range("AN23").Value = 5
$ws.Range("AO23").Value = 251
$ws.Range("AP23").Value = 3
$ws.Range("AQ23").Value = 2
$ws.Range("AR23").Value = 1
$ws.Range("AS23").Value = 6
$ws.Range("AT23").Value = 2
$ws.Range("AU23").Value = 3
$ws.Range("AW23").Value = 2
$ws.Range("AX23").Value = 3
$ws.Range("AY23").Value = 6
$ws.Range("BA23").Value = 5

# --- Restore the selection to match the authored view state. ---
$ws.Range("BB23").Select() | Out-Null
